$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos (author: "fixed typos in elemente excel")
$ws.Range("B15").Value = "Als Waffe werden in der Regel alle Gegenstände bezeichnet, die Fähig sind, Lebewesen Schaden zuzufügen."
$ws.Range("B23").Value = "Du… und andere, die so sind wie du sind…"
$ws.Range("B25").Value = "Bäume bestehen aus Papier und wachsen im Baumarkt."

# Reflect the author's final cursor position/selection in the sheet view
$ws.Range("B15").Select()
